$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for rows with changed values.
# D-column values are plain text in the source data (e.g. "29.364.39", "0.7130"),
# so we force a Text number format before assignment and then clear the format
# again so the stored cell keeps its original (default) style while remaining a
# text value instead of being auto-coerced into a number (which would silently
# drop meaningful trailing zeros / group separators).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.364.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.89"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7130"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08094"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3135"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08348"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.20"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7215"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.250"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.00"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.290"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.93%  "

$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.371.95"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.126.65"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.806"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.080"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.427"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.348"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.218"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05384"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.956"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.78%  "

# Rows 35 and 36 swap coin identity (ImmutableX <-> ARBITRUM) with updated values
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.180"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7521"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01882"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.283.84"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.746"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.572"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.18%  "

$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8931"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "110.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.26%  "

$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("E46").Value = "  +6.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.012.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5213"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.488"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4367"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.48%  "
